$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (GUI/Excel formatting tweak)
$ws.Columns.Item(8).ColumnWidth = 17.833333333333332
$ws.Columns.Item(10).ColumnWidth = 18.833333333333332
$ws.Columns.Item(17).ColumnWidth = 19.833333333333332
$ws.Columns.Item(21).ColumnWidth = 18.833333333333332

# Updated error-calculation results for row 2
$ws.Range("H2").Value = 3.125603062778264
$ws.Range("I2").Value = 0.01283955242086306
$ws.Range("J2").Value = 3.54551924390971
$ws.Range("K2").Value = 1.069781537870688
$ws.Range("P2").Value = 0.5832367987382572
$ws.Range("Q2").Value = 0.3212830336741677
$ws.Range("T2").Value = 0.002070408906527906
$ws.Range("U2").Value = 0.2893874872767901
